# Updates for figures 2 and 5 following editorial requests.
#
# Figure 2 (sheet "top10ledge"): the Orycteropus afer (aardvark) photo
# credit/source is swapped from David Renoult's iNaturalist photo to
# Dave Brown's.
#
# Figure 5 (sheet "introduced"): no data change, but the active/selected
# worksheet moves from "introduced" to "top10ledge".

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("top10ledge")

# Row 2 (Orycteropus afer): replace the photo URL + credit.
$ws2.Range("B2").Value = "https://www.inaturalist.org/photos/78448701"
$ws2.Range("C2").Value = "Dave Brown"

# The url cell becomes a hyperlink (matching the other url_photo cells),
# styled with the workbook's built-in hyperlink style.
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://www.inaturalist.org/photos/78448701")
$ws2.Range("B2").Style = "Lien hypertexte"

# Editorial re-review moved the active tab/selection from "introduced" to
# "top10ledge".
$ws2.Activate()
$ws2.Range("B15").Select()
